# Insert two new weekly price observations into the "Hortaliza, Feria
# Lagunitas de Puerto Montt - Cebollín" consolidated sheet.
#
# The sheet is an append-log of weekly price records (one per market /
# category / quality / week) that is not chronologically sorted. This
# edit adds one new observation near the top of the historical block
# (landing at row 368) and a second new observation further down
# (landing at row 481, using the row numbering that results after the
# first insert), pushing every following row down by one each time the
# dimension grows from A1:R483 to A1:R485.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-Row($Sheet, $Row, $Fecha, $Vol, $Min, $Max, $Avg, $Kg) {
    $Sheet.Range("A$Row").Value = 4
    $Sheet.Range("B$Row").Value = "Feria Lagunitas de Puerto Montt"
    $Sheet.Range("C$Row").Value = "Los Lagos"
    $Sheet.Range("D$Row").Value = $Fecha
    $Sheet.Range("E$Row").Value = 10
    $Sheet.Range("F$Row").Value = 100112037
    $Sheet.Range("G$Row").Value = "Cebollín"
    $Sheet.Range("H$Row").Value = "Sin especificar"
    $Sheet.Range("I$Row").Value = "Primera"
    $Sheet.Range("J$Row").Value = $Vol
    $Sheet.Range("K$Row").Value = $Min
    $Sheet.Range("L$Row").Value = $Max
    $Sheet.Range("M$Row").Value = $Avg
    $Sheet.Range("N$Row").Value = "`$/paquete 36 unidades"
    $Sheet.Range("O$Row").Value = "Región Metropolitana"
    $Sheet.Range("P$Row").Value = $Kg
    $Sheet.Range("Q$Row").Value = 36
    $Sheet.Range("R$Row").Value = "Hortaliza"
}

# First insertion: a new row at (current) row 368 — shifts old rows
# 368..483 down to 369..484.
$ws.Rows(368).Insert()
Set-Row $ws 368 45120 70 6500 6500 6500 181

# Second insertion: a new row at (current) row 481 — shifts what is now
# rows 481..484 down to 482..485.
$ws.Rows(481).Insert()
Set-Row $ws 481 45121 160 6500 7500 7000 194
